$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Content fix: "nm_4400_fix13.sql" -> "nm_4500_fix13.sql" and the
#    adjoining version cell "3.0" -> "1.0" in the file-list table.
# ------------------------------------------------------------------
$fileTable = $d.Tables.Item(2)
for ($r = 1; $r -le $fileTable.Rows.Count; $r++) {
    $cellText = $fileTable.Cell($r, 1).Range.Text
    if ($cellText -like "*nm_4400_fix13.sql*") {
        $nameRange = $fileTable.Cell($r, 1).Range
        $nameRange.MoveEnd(1, -1) | Out-Null   # drop the end-of-cell marker
        $nameRange.Find.Execute("nm_4400_fix13.sql", $true, $false, $false, $false, $false, $true, 0, $false, "nm_4500_fix13.sql", 1) | Out-Null

        $verRange = $fileTable.Cell($r, 2).Range
        $verRange.MoveEnd(1, -1) | Out-Null
        $verRange.Find.Execute("3.0", $true, $false, $false, $false, $false, $true, 0, $false, "1.0", 1) | Out-Null
    }
}

# ------------------------------------------------------------------
# 2) Bookmark housekeeping: Word re-drops the "_GoBack" bookmark at the
#    site of the last edit (an empty paragraph near the top of the
#    document) every time the file is saved after an edit, moving it
#    away from wherever it previously sat (the final paragraph here).
#    Adding a bookmark with the reserved name "_GoBack" automatically
#    relocates/removes any prior one and renumbers the remaining
#    bookmarks, matching the diff.
# ------------------------------------------------------------------
$target = $d.Content.Paragraphs.Item(7).Range
$d.Bookmarks.Add("_GoBack", $target) | Out-Null
